$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("19Tto20TMap")

# Increment the period (column I, "pe") for every data row from 2019Oct to 2020Oct
for ($r = 2; $r -le 91; $r++) {
    $cell = $ws.Cells.Item($r, 9)
    if ($cell.Value -eq "2019Oct") {
        $cell.Value = "2020Oct"
    }
}

# Restore the view/selection state to match target
$ws.Activate()
$ws.Application.ActiveWindow.SplitColumn = 1
$ws.Application.ActiveWindow.FreezePanes = $true

$ws.Range("I1:I1048576").Select()
$ws.Range("A44").Select()
